$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates scraped from the commit diff: columns B (Coin), C (Link),
# D (Price) and E (Volume 1h) are refreshed for rows 2-51. Column D holds
# text that often *looks* numeric (e.g. "0.460", "603.21"); a leading
# apostrophe forces Excel to store it as text (matching the original
# inlineStr cells) instead of silently converting it to a float, and the
# follow-up Style reset drops the auto-applied quote-prefix style so the
# cell stays unstyled like the source file.
$updates = @(
    @{Row=2; D="68.820.87"; E="  +0.10%  "},
    @{Row=3; D="3.855.17"; E="  -0.26%  "},
    @{Row=4; D="0.998"; E="  -0.23%  "},
    @{Row=5; D="603.21"; E="  +0.20%  "},
    @{Row=6; D="169.18"; E="  +0.40%  "},
    @{Row=7; D="3.858.73"; E="  -0.24%  "},
    @{Row=8; E="  -0.03%  "},
    @{Row=9; D="0.527"; E="  -0.63%  "},
    @{Row=10; E="  +0.64%  "},
    @{Row=11; D="6.48"; E="  +2.03%  "},
    @{Row=12; D="0.460"; E="  -1.08%  "},
    @{Row=13; D="0.0000278"; E="  +11.01%  "},
    @{Row=14; D="36.93"; E="  -2.07%  "},
    @{Row=15; D="4.498.90"; E="  -0.25%  "},
    @{Row=16; D="3.845.84"; E="  -0.45%  "},
    @{Row=17; D="68.738.39"; E="  -0.18%  "},
    @{Row=18; D="18.36"; E="  -0.92%  "},
    @{Row=19; D="7.35"; E="  -2.95%  "},
    @{Row=20; E="  -0.60%  "},
    @{Row=21; D="10.99"; E="  +1.45%  "},
    @{Row=22; D="472.29"; E="  -1.62%  "},
    @{Row=23; D="0.726"; E="  -1.66%  "},
    @{Row=24; D="0.0000163"; E="  +0.28%  "},
    @{Row=25; D="83.37"; E="  -1.50%  "},
    @{Row=26; E="  -0.49%  "},
    @{Row=27; D="12.08"; E="  -2.08%  "},
    @{Row=28; D="10.42"; E="  +3.28%  "},
    @{Row=29; E="  +0.14%  "},
    @{Row=30; D="2.96"},
    @{Row=31; D="4.001.71"; E="  -0.27%  "},
    @{Row=32; D="7.69"; E="  -0.92%  "},
    @{Row=33; D="31.33"; E="  +0.10%  "},
    @{Row=34; E="  -1.05%  "},
    @{Row=35; D="9.28"; E="  -2.50%  "},
    @{Row=36; D="3.820.10"; E="  -0.29%  "},
    @{Row=37; B="dogwifhat"; C="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D="3.76"; E="  +14.81%  "},
    @{Row=38; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.104"; E="  -1.47%  "},
    @{Row=39; E="  +0.17%  "},
    @{Row=40; D="0.140"; E="  -0.07%  "},
    @{Row=41; D="5.92"},
    @{Row=42; D="0.997"; E="  -0.30%  "},
    @{Row=43; E="  -0.37%  "},
    @{Row=44; E="  -1.63%  "},
    @{Row=45; B="USDe"; C="https://coinranking.com/coin/exbfr2U-0+usde-usde"; D="1.00"; E="  -0.02%  "},
    @{Row=46; B="Bittensor"; C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D="421.08"; E="  -2.04%  "},
    @{Row=47; B="FLOKI"; C="https://coinranking.com/coin/fmHk13Rqw+floki-floki"; D="0.000296"; E="  +10.07%  "},
    @{Row=48; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="8.67"; E="  +0.41%  "},
    @{Row=49; D="46.97"; E="  -1.80%  "},
    @{Row=50; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="26.23"; E="  +5.75%  "},
    @{Row=51; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="141.77"; E="  -0.12%  "},
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
        $ws.Cells.Item($u.Row, 4).Style = "Normal"
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
